$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 234.72728
$ws.Range("I4").Value = 144.66667
$ws.Range("J4").Value = 640
$ws.Range("K4").Value = 144.66667
$ws.Range("L4").Value = 640
$ws.Range("M4").Value = -30.66667000000001
$ws.Range("N4").Value = -868
$ws.Range("H19").Value = 2700
$ws.Range("I19").Value = 750
$ws.Range("J19").Value = 4000
$ws.Range("K19").Value = 750
$ws.Range("L19").Value = 4000
$ws.Range("M19").Value = -575
$ws.Range("N19").Value = -4350
$ws.Range("H43").Value = 9276893
$ws.Range("J43").Value = 13890089
$ws.Range("L43").Value = 13890089
$ws.Range("N43").Value = -13890227
$ws.Range("H44").Value = 23000
$ws.Range("J44").Value = 23000
$ws.Range("L44").Value = 23000
$ws.Range("N44").Value = -23924
$ws.Range("H64").Value = 3922.7273
$ws.Range("I64").Value = 4338
$ws.Range("J64").Value = 3576.6667
$ws.Range("K64").Value = 4338
$ws.Range("L64").Value = 3576.6667
$ws.Range("M64").Value = -4090
$ws.Range("N64").Value = -4072.6667
$ws.Range("H67").Value = 3922.7273
$ws.Range("I67").Value = 4338
$ws.Range("J67").Value = 3576.6667
$ws.Range("K67").Value = 4338
$ws.Range("L67").Value = 3576.6667
$ws.Range("M67").Value = -3480
$ws.Range("N67").Value = -5292.6667
$ws.Range("H88").Value = 951702.7
$ws.Range("I88").Value = 995
$ws.Range("J88").Value = 1030928.3
$ws.Range("K88").Value = 995
$ws.Range("L88").Value = 1030928.3
$ws.Range("M88").Value = -589
$ws.Range("N88").Value = -1031740.3
$ws.Range("H91").Value = 951702.7
$ws.Range("I91").Value = 995
$ws.Range("J91").Value = 1030928.3
$ws.Range("K91").Value = 995
$ws.Range("L91").Value = 1030928.3
$ws.Range("M91").Value = 409
$ws.Range("N91").Value = -1033736.3
$ws.Range("H94").Value = 5333.3335
$ws.Range("I94").Value = 6500
$ws.Range("K94").Value = 6500
$ws.Range("M94").Value = -6049
$ws.Range("H116").Value = 3452.4285
$ws.Range("I116").Value = 3008.75
$ws.Range("J116").Value = 4044
$ws.Range("K116").Value = 3008.75
$ws.Range("L116").Value = 4044
$ws.Range("M116").Value = 433.25
$ws.Range("N116").Value = -10928
$ws.Range("H132").Value = 10757669
$ws.Range("I132").Value = 11907448
$ws.Range("J132").Value = 26400
$ws.Range("K132").Value = 35722344
$ws.Range("L132").Value = 79200
$ws.Range("M132").Value = -35719814
$ws.Range("N132").Value = -84260
$ws.Range("H135").Value = 830.3
$ws.Range("I135").Value = 305.46155
$ws.Range("J135").Value = 1805
$ws.Range("K135").Value = 2749.15395
$ws.Range("L135").Value = 16245
$ws.Range("M135").Value = -214.1539499999999
$ws.Range("N135").Value = -21315
$ws.Range("H138").Value = 1195.914
$ws.Range("I138").Value = 773.3461
$ws.Range("J138").Value = 1731.8536
$ws.Range("K138").Value = 2320.0383
$ws.Range("L138").Value = 5195.560799999999
$ws.Range("M138").Value = 2819.9617
$ws.Range("N138").Value = -15475.5608
$ws.Range("H141").Value = 1193.4286
$ws.Range("J141").Value = 3000
$ws.Range("L141").Value = 9000
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2391.1428
$ws.Range("I88").Value = 1663.75
$ws.Range("J88").Value = 2838.7693
$ws.Range("K88").Value = 1663.75
$ws.Range("L88").Value = 2838.7693
$ws.Range("M88").Value = -1257.75
$ws.Range("N88").Value = -3650.7693
$ws.Range("H91").Value = 2391.1428
$ws.Range("I91").Value = 1663.75
$ws.Range("J91").Value = 2838.7693
$ws.Range("K91").Value = 1663.75
$ws.Range("L91").Value = 2838.7693
$ws.Range("M91").Value = -259.75
$ws.Range("N91").Value = -5646.7693
$ws.Range("H97").Value = 551
$ws.Range("I97").Value = 400
$ws.Range("J97").Value = 1004
$ws.Range("K97").Value = 400
$ws.Range("L97").Value = 1004
$ws.Range("M97").Value = 96
$ws.Range("N97").Value = -1996
$ws.Range("H102").Value = 166666670
$ws.Range("I102").Value = 166666670
$ws.Range("K102").Value = 166666670
$ws.Range("M102").Value = -166665048
$ws.Range("H103").Value = 75000
$ws.Range("J103").Value = 75000
$ws.Range("L103").Value = 75000
$ws.Range("N103").Value = -77344
$ws.Range("H132").Value = 2725.4
$ws.Range("I132").Value = 1939.5555
$ws.Range("J132").Value = 3904.1667
$ws.Range("K132").Value = 5818.666499999999
$ws.Range("L132").Value = 11712.5001
$ws.Range("M132").Value = -3288.666499999999
$ws.Range("N132").Value = -16772.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3927
$ws.Range("I86").Value = 4323.2104
$ws.Range("K86").Value = 4323.2104
$ws.Range("M86").Value = -3200.2104
$ws.Range("H89").Value = 3927
$ws.Range("I89").Value = 4323.2104
$ws.Range("K89").Value = 21616.052
$ws.Range("M89").Value = -16000.052
$ws.Range("H99").Value = 38462940
$ws.Range("I99").Value = 41668020
$ws.Range("K99").Value = 41668020
$ws.Range("M99").Value = -41666522
$ws.Range("H105").Value = 100002536
$ws.Range("J105").Value = 2574.75
$ws.Range("L105").Value = 2574.75
$ws.Range("N105").Value = -6068.75
$ws.Range("H134").Value = 8771.875
$ws.Range("I134").Value = 841
$ws.Range("K134").Value = 2523
$ws.Range("M134").Value = 12

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1335.7142
$ws.Range("I31").Value = 1316.6666
$ws.Range("J31").Value = 1450
$ws.Range("K31").Value = 1316.6666
$ws.Range("L31").Value = 1450
$ws.Range("M31").Value = -1021.6666
$ws.Range("N31").Value = -2040
$ws.Range("H34").Value = 1335.7142
$ws.Range("I34").Value = 1316.6666
$ws.Range("J34").Value = 1450
$ws.Range("K34").Value = 1316.6666
$ws.Range("L34").Value = 1450
$ws.Range("M34").Value = -1114.6666
$ws.Range("N34").Value = -1854
$ws.Range("H99").Value = 1774.5385
$ws.Range("I99").Value = 1824.4546
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1824.4546
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -326.4546
$ws.Range("N99").Value = -4496
$ws.Range("H105").Value = 510.33334
$ws.Range("H107").Value = 706.9167
$ws.Range("J107").Value = 914.6667
$ws.Range("L107").Value = 914.6667
$ws.Range("N107").Value = -4754.6667
$ws.Range("H109").Value = 29714.428
$ws.Range("J109").Value = 29714.428
$ws.Range("L109").Value = 29714.428
$ws.Range("N109").Value = -31794.428
$ws.Range("H126").Value = 1774.5385
$ws.Range("I126").Value = 1824.4546
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 5473.3638
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -3003.3638
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 8824.666999999999
$ws.Range("I132").Value = 15996
$ws.Range("J132").Value = 4261.091
$ws.Range("K132").Value = 47988
$ws.Range("L132").Value = 12783.273
$ws.Range("M132").Value = -45458
$ws.Range("N132").Value = -17843.273
$ws.Range("H134").Value = 23812012
$ws.Range("I134").Value = 30305562
$ws.Range("K134").Value = 90916686
$ws.Range("M134").Value = -90914151

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2249.7727
$ws.Range("I102").Value = 2416.389
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 2416.389
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -794.3890000000001
$ws.Range("N102").Value = -4744
$ws.Range("H126").Value = 1965.4667
$ws.Range("I126").Value = 1798.5834
$ws.Range("K126").Value = 5395.7502
$ws.Range("M126").Value = -2925.7502
$ws.Range("H132").Value = 2256
$ws.Range("I132").Value = 1918.5652
$ws.Range("J132").Value = 3226.125
$ws.Range("K132").Value = 5755.6956
$ws.Range("L132").Value = 9678.375
$ws.Range("M132").Value = -3225.6956
$ws.Range("N132").Value = -14738.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 659.8
$ws.Range("I16").Value = 750.7692
$ws.Range("J16").Value = 490.85715
$ws.Range("K16").Value = 750.7692
$ws.Range("L16").Value = 490.85715
$ws.Range("M16").Value = -580.7692
$ws.Range("N16").Value = -830.85715
$ws.Range("H22").Value = 1723.1428
$ws.Range("I22").Value = 1380
$ws.Range("J22").Value = 1780.3334
$ws.Range("K22").Value = 1380
$ws.Range("L22").Value = 1780.3334
$ws.Range("M22").Value = -1085
$ws.Range("N22").Value = -2370.3334
$ws.Range("H27").Value = 1723.1428
$ws.Range("I27").Value = 1380
$ws.Range("J27").Value = 1780.3334
$ws.Range("K27").Value = 1380
$ws.Range("L27").Value = 1780.3334
$ws.Range("M27").Value = -1273
$ws.Range("N27").Value = -1994.3334
$ws.Range("H132").Value = 24502.318
$ws.Range("I132").Value = 1438.4783
$ws.Range("J132").Value = 49762.715
$ws.Range("K132").Value = 4315.4349
$ws.Range("L132").Value = 149288.145
$ws.Range("M132").Value = -1785.4349
$ws.Range("N132").Value = -154348.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 680
$ws.Range("I81").Value = 680
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1360
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -299
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 680
$ws.Range("I84").Value = 680
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 6800
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -1496
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 13687212
$ws.Range("I122").Value = 16253201
$ws.Range("K122").Value = 48759603
$ws.Range("M122").Value = -48757153
$ws.Range("H132").Value = 3639.7222
$ws.Range("I132").Value = 3656.0908
$ws.Range("J132").Value = 3614
$ws.Range("K132").Value = 10968.2724
$ws.Range("L132").Value = 10842
$ws.Range("M132").Value = -8438.2724
$ws.Range("N132").Value = -15902
